$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "XD86" -> "CD86" in the Monocytes Macrophage marker gene list (row 10, col C)
$ws.Range("C10").Value = "CD163, CD86, CD80, CD68, CD74, CD14"

# Move the active selection to C10 (the cell that was edited), matching the
# author moving focus there after correcting the gene list / enabling
# selection of the default assay for Module Scoring.
$ws.Range("C10").Select()
